# Update the price list dated in cell A1 (month rolled forward) and the
# unit-price column (D) for rows 14-21 on the "PLANCHUELA DISMAY" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date in A1: 2024-04-24 (serial 45406) -> 2024-05-24 (serial 45436)
$ws.Range("A1").Value = 45436

# Updated unit prices (column D) for rows 14-21
$ws.Range("D14").Value = 83.175
$ws.Range("D15").Value = 108.235
$ws.Range("D16").Value = 129.622
$ws.Range("D17").Value = 207.394
$ws.Range("D18").Value = 259.243
$ws.Range("D19").Value = 324.054
$ws.Range("D20").Value = 360.78
$ws.Range("D21").Value = 399.666
